$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "PerformanceEvaluationId"
$ws.Range("B1").Value = "EmployeeId"
$ws.Range("C1").Value = "EvaluationDate"
$ws.Range("D1").Value = "ConsciousnessScore"
$ws.Range("E1").Value = "DisciplinaryViolations"
$ws.Range("F1").Value = "UnexcusedAbsences"
$ws.Range("G1").Value = "WorkPerformanceScore"
$ws.Range("H1").Value = "OvertimeHours"

$ws.Range("L16").Select()
